$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  "05-050305TC" = "a0c1161837786ed577bd398f7504ad26"
  "05-050301A" = "5e3060384c34d729b00004f7ae145871"
  "05-050207TP" = "f70941508953fb002528f6913831f208"
  "05-050305TP" = "dbb17ca4b52a4c7e5e94472e9b66584d"
  "05-050316TC" = "b5f6d0d26190225d3934df2a1696f2eb"
  "05-050302A" = "40e9a17fbb91f23587413a1d9790a202"
  "05-050316TP" = "e0c01ed14f14658dee48c59beec4108d"
  "05-050301TP" = "89e31980121a03ecb4d90a72f238e8a1"
  "05-050312TP" = "d5cb4a103f42f397f487757d5d1944b5"
  "05-050312TC" = "fdf78d44d331f49c30cdad3e5f603af7"
  "05-050308A" = "f26197f222ebf1ddb0efdcaf398412ff"
  "05-050303TC" = "5d06c54d6ae33de6ef32267596888733"
  "05-050305A" = "bfd43c2f789ae217aee9d6a0c58b3db0"
  "05-050314TP" = "4d6ab91e2c46180e790ad8b177c98dfb"
  "05-050314TC" = "e092507be44a5fc8a934d3c321d61312"
  "05-050306A" = "d390d1e05d7bb974a4a42141255c0a4f"
  "05-050303A" = "6237ab96bee52eeaaab382477910cf24"
  "05-050205TP" = "38806a07acc53d2c72acefaafded9c2a"
  "05-050304A" = "32971f9e01b7e44aa184d2c517c589a3"
  "05-050205TC" = "a160247022553052b8abae884d3456af"
  "05-050313A" = "c52c7a90e94c5465e55e2c08fa133e27"
  "05-050205A" = "c13f2fbb39889312ebe51655c8906f50"
  "05-050314A" = "466dbbdd75b29f19f4337f486efc3f50"
  "05-050311A" = "0609d0fedd5ea75cbda7ed438647537d"
  "05-050306TP" = "c3a5a348978c43abef5433eb005a5bcf"
  "05-050317TC" = "ddb1b7c47afac42d25035f85f87784b0"
  "05-050317TP" = "f233a859f5bd79f16e9f2757fc8bd8f3"
  "05-050308TP" = "c3f60ea1fa19ab1c30e5690afe2c4a50"
  "05-050302TP" = "d800fd7c62ec6c0c850576a265f8b098"
  "05-050313TP" = "54a38eade1ac1d27fbf81691c24ba5cb"
  "05-050313TC" = "a5cbd6b59e2b8d03a54f6e0298e7be71"
  "05-050317A" = "ad397fcbf26ca599eacf01feded20657"
  "05-050206TP" = "be0392d0a6d60ca7e9618cc528ca05c4"
  "05-050304TC" = "f50477f816aa871fa3ad5dcebd3db88f"
  "05-050206A" = "a0988d882b880d8a537e3239a8df8a88"
  "05-050315A" = "b384debc4d274cd463c47814bb03584f"
  "05-050304TP" = "7068eaeabb596cde9800331635f8126a"
  "05-050316A" = "facd5956c8a905ffec71ab1f34091b5f"
  "05-050207A" = "46fd746591fd901150cf8faf3cfdf069"
  "05-050315TP" = "a8bb0f961ff788315d8ff7476ab13e9e"
  "05-050311TC" = "8e6ce109c0d4d26ba84457b761a405fb"
  "05-050311TP" = "f8983d626543a44eea837a380749594b"
  "05-050309TC" = "489eb9999e6bb05525b849268980c21f"
}

foreach ($key in $updates.Keys) {
  $found = $ws.Columns(1).Find($key, [Type]::Missing, [Type]::Missing, 1)
  if ($found -ne $null) {
    $row = $found.Row
    $ws.Cells.Item($row, 2).Value = $updates[$key]
  } else {
    Write-Output ("Key not found: " + $key)
  }
}

Write-Output "Done"
